# Update Handback report timestamps (Generate Report for Handback)
$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 01:08:58"
$wsZhCn.Range("H2").Value = "2016-03-23 01:09:19"

# de-de sheet: Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 01:09:02"
$wsDeDe.Range("H2").Value = "2016-03-23 01:09:26"
